$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new customer record row (row 8): name "Doe" and its image path.
$ws.Range("A8").Value = "Doe"
$ws.Range("B8").Value = "D:\Robot\myproject\hello\images/p2.JPG"
